# CI: Update Excel counters (state_counters + packages)
#
# Appends 15 new package rows (rows 46-60) to the "Packages" worksheet,
# extending the used range from A1:G45 to A1:G60. Columns are:
#   A = PackageType, B = State, C = Authority, D = ActionType,
#   E = PackageID,   F = Status, G = ParentID
#
# Some rows have blank ActionType/ParentID cells (SPA-style rows). A plain
# empty-string assignment clears/omits the cell entirely in this engine, so
# for those we write a leading apostrophe (forces a text cell containing an
# empty string, reusing the workbook's existing blank shared string) and
# then immediately reset the cell style to "Normal" so it does not keep the
# implicit quote-prefix formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    if ($value -eq "") {
        $ws.Range($addr).Value = "'"
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $value
    }
}

$rows = @(
    @(46, "Waiver", "MD", "1915(c)",     "Amendment", "MD-2260.R00.62", "Submitted",           "MD-2260.R00.00"),
    @(47, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9531",     "Under Review",        ""),
    @(48, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9532",     "Under Review",        ""),
    @(49, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9533",     "Disapproved",         ""),
    @(50, "Waiver", "MD", "1915(c)",     "Amendment", "MD-2260.R00.63", "Pending-Approval",    "MD-2260.R00.00"),
    @(51, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9534",     "Pending-Concurrence", ""),
    @(52, "SPA",    "MD", "CHIP SPA",    "",          "MD-25-9535",     "Submitted",           ""),
    @(53, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9536",     "RAI Issued",          ""),
    @(54, "SPA",    "MD", "CHIP SPA",    "",          "MD-25-9537",     "Submitted",           ""),
    @(55, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9538",     "Submitted",           ""),
    @(56, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9539",     "Submitted",           ""),
    @(57, "Waiver", "MD", "1915(c)",     "Amendment", "MD-2260.R00.64", "Unsubmitted",         "MD-2260.R00.00"),
    @(58, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9540",     "Under Review",        ""),
    @(59, "Waiver", "MD", "1915(b)",     "Initial",   "MD-2284.R00.00", "Terminated",          ""),
    @(60, "SPA",    "MD", "Medicaid SPA","",          "MD-25-9541",     "Withdrawn",           "")
)

foreach ($row in $rows) {
    $r = $row[0]
    Set-Text ("A" + $r) $row[1]
    Set-Text ("B" + $r) $row[2]
    Set-Text ("C" + $r) $row[3]
    Set-Text ("D" + $r) $row[4]
    Set-Text ("E" + $r) $row[5]
    Set-Text ("F" + $r) $row[6]
    Set-Text ("G" + $r) $row[7]
}

Write-Output "Added rows 46-60 to Packages sheet"
